$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# D-column values are plain digit/dot strings that Excel would otherwise
# auto-coerce to Number (losing formatting like trailing zeros / multi-dot
# grouping), so force Text format before assigning them.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.037.20'
$ws.Range('E2').Value = '  -0.14%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.650.32'
$ws.Range('E3').Value = '  +0.30%  '

$ws.Range('E4').Value = '  -0.31%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.36'
$ws.Range('E5').Value = '  +0.40%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5199'
$ws.Range('E6').Value = '  +0.38%  '

$ws.Range('E7').Value = '  -0.29%  '

$ws.Range('E8').Value = '  +0.95%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06325'
$ws.Range('E9').Value = '  +0.92%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.40'
$ws.Range('E10').Value = '  +0.46%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07655'
$ws.Range('E11').Value = '  -1.44%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.584'
$ws.Range('E12').Value = '  +2.75%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.681.79'
$ws.Range('E13').Value = '  +1.36%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.877.27'
$ws.Range('E14').Value = '  +0.13%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5602'
$ws.Range('E15').Value = '  +1.37%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅8151'
$ws.Range('E16').Value = '  +2.48%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.23'
$ws.Range('E17').Value = '  +0.93%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '26.037.74'
$ws.Range('E18').Value = '  -0.16%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.616'
$ws.Range('E20').Value = '  +0.14%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.50'
$ws.Range('E21').Value = '  +4.52%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '191.62'
$ws.Range('E22').Value = '  -0.74%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.905'
$ws.Range('E23').Value = '  -0.43%  '

$ws.Range('E24').Value = '  -0.38%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.57'
$ws.Range('E25').Value = '  -2.30%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1185'
$ws.Range('E26').Value = '  -1.20%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.195'
$ws.Range('E27').Value = '  +0.76%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.86'
$ws.Range('E28').Value = '  -0.13%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.510'
$ws.Range('E29').Value = '  +2.41%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05396'
$ws.Range('E30').Value = '  -3.74%  '

$ws.Range('E31').Value = '  +0.34%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.447'
$ws.Range('E32').Value = '  -0.90%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.352'
$ws.Range('E33').Value = '  -0.18%  '

$ws.Range('E34').Value = '  -2.11%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.424'
$ws.Range('E35').Value = '  +0.81%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.780'
$ws.Range('E36').Value = '  -0.71%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9451'
$ws.Range('E37').Value = '  +1.35%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5629'
$ws.Range('E38').Value = '  -0.15%  '

$ws.Range('E39').Value = '  +0.20%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.863'
$ws.Range('E40').Value = '  -1.03%  '

$ws.Range('E41').Value = '  -0.20%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.028.96'
$ws.Range('E42').Value = '  -2.85%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8249'
$ws.Range('E43').Value = '  -1.30%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.08'
$ws.Range('E44').Value = '  -1.56%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.787.42'
$ws.Range('E45').Value = '  +0.06%  '

$ws.Range('E46').Value = '  +4.54%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '57.39'
$ws.Range('E47').Value = '  +1.02%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.9991'
$ws.Range('E48').Value = '  -0.18%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4320'
$ws.Range('E49').Value = '  -0.32%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.941'
$ws.Range('E50').Value = '  +0.02%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05143'
$ws.Range('E51').Value = '  -3.43%  '
